# Atualizando o arquivo XLSX
# Apply updated odds values per row, matching the upstream FlashScore refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("G7").Value = 1.42
$ws.Range("H7").Value = 4.33
$ws.Range("I7").Value = 8.5
$ws.Range("J7").Value = 1.95
$ws.Range("K7").Value = 2.2
$ws.Range("L7").Value = 8.5
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.67
$ws.Range("U7").Value = 2.38
$ws.Range("V7").Value = 1.53
$ws.Range("X7").Value = 5.5
$ws.Range("Z7").Value = 8.5
$ws.Range("AA7").Value = 13
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 8.5
$ws.Range("AE7").Value = 26
$ws.Range("AF7").Value = 101
$ws.Range("AH7").Value = 15
$ws.Range("AI7").Value = 41
$ws.Range("AJ7").Value = 26
$ws.Range("AK7").Value = 101
$ws.Range("AL7").Value = 67
$ws.Range("AN7").Value = 3.1
$ws.Range("AO7").Value = 7
$ws.Range("AQ7").Value = 21
$ws.Range("AU7").Value = 11
$ws.Range("AW7").Value = 126
$ws.Range("AX7").Value = 9
$ws.Range("AZ7").Value = 51
$ws.Range("BA7").Value = 251
$ws.Range("BB7").Value = 251

# Row 8
$ws.Range("G8").Value = 1.9
$ws.Range("I8").Value = 4.5
$ws.Range("J8").Value = 2.63
$ws.Range("Q8").Value = 2.4
$ws.Range("U8").Value = 2.2
$ws.Range("V8").Value = 1.62
$ws.Range("X8").Value = 7.5
$ws.Range("Z8").Value = 15
$ws.Range("AE8").Value = 21
$ws.Range("AN8").Value = 3.6
$ws.Range("AY8").Value = 29

# Row 9
$ws.Range("G9").Value = 2.57
$ws.Range("I9").Value = 3
$ws.Range("K9").Value = 1.77
$ws.Range("M9").Value = 1.14
$ws.Range("N9").Value = 5.5
$ws.Range("Q9").Value = 3.1
$ws.Range("R9").Value = 1.36
$ws.Range("S9").Value = 1.73
$ws.Range("T9").Value = 2
$ws.Range("V9").Value = 1.5
$ws.Range("AE9").Value = 21
$ws.Range("AU9").Value = 10
$ws.Range("BA9").Value = 67
$ws.Range("BB9").Value = 126

# Row 15
$ws.Range("G15").Value = 1.5
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 6.5
$ws.Range("J15").Value = 2.05
$ws.Range("K15").Value = 2.37
$ws.Range("L15").Value = 6
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 13
$ws.Range("O15").Value = 1.22
$ws.Range("P15").Value = 4
$ws.Range("U15").Value = 1.83
$ws.Range("V15").Value = 1.83
$ws.Range("W15").Value = 7
$ws.Range("AH15").Value = 17
$ws.Range("AI15").Value = 34
$ws.Range("AJ15").Value = 21
$ws.Range("AK15").Value = 67
$ws.Range("AO15").Value = 7.5
$ws.Range("AS15").Value = 126
$ws.Range("AX15").Value = 7.5
$ws.Range("BC15").Value = 251

# Row 17
$ws.Range("O17").Value = 1.3
$ws.Range("P17").Value = 3.4
$ws.Range("Q17").Value = 1.98
$ws.Range("R17").Value = 1.83
$ws.Range("U17").Value = 1.73
$ws.Range("V17").Value = 2
$ws.Range("AB17").Value = 26
$ws.Range("AG17").Value = 201
$ws.Range("BC17").Value = 151

# Row 47
$ws.Range("N47").Value = 13
$ws.Range("T47").Value = 3.54
